# Debugging Marksheet - fill in student results, compute totals, and rank them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("K1").Value = "TOTAL"

# ---- Student rows (pre-sort order; column A already holds Sr No 1..9) ----
# Row layout per student: Name, then marks for Code_0..Code_6 (columns D..J)
$names = @(
    "John Sijo",
    "Hanushree M",
    "Sidharth KS",
    "Dhamodharan S P",
    "Arshin Joseph Giril ",
    "Aswin ES",
    "Aahil Muhammed",
    "Alwin Joshy",
    "Naveenkumar T"
)

$marks = @(
    @(0, 0, 0, 0, 0, 10, 20),
    @(0, 0, 0, 0, 0, 0, 12),
    @(10, 0, 10, 0, 10, 0, 16),
    @(10, 0, 10, 10, 10, 10, 40),
    @(0, 0, 0, 0, 0, 0, 20),
    @(10, 0, 0, 0, 0, 0, 12),
    @(10, 10, 10, 10, 10, 0, 24),
    @(0, 0, 0, 0, 0, 0, 12),
    @(0, 10, 0, 0, 10, 0, 12)
)

# One row (Dhamodharan) was entered with a literal total instead of a formula.
$literalTotalRow = 3

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]

    $col = 4
    foreach ($mark in $marks[$i]) {
        $ws.Cells.Item($row, $col).Value = $mark
        $col = $col + 1
    }

    if ($i -eq $literalTotalRow) {
        $ws.Cells.Item($row, 11).Value = 80
    } else {
        $ws.Cells.Item($row, 11).Formula = "=(D$row+E$row+F$row+G$row+H$row+I$row+J$row)"
    }
}

# ---- Remove the old trailing placeholder row ----
$ws.Rows("11").Delete()

# ---- Sort students by TOTAL (column K) descending ----
$sf = $ws.Sort.SortFields
$sf.Add($ws.Range("K2:K10"), 0, 2)
$ws.Sort.SetRange($ws.Range("A1:L10"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# ---- Rank column, assigned top-to-bottom after sorting ----
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 12).Value = $r - 1
}

# ---- View state ----
$win = $ws.Application.ActiveWindow
$win.Zoom = 180
$ws.Range("L12").Select() | Out-Null
